# "split methods by part" - push the "Testing" and "Postmortem" log rows
# down by two rows (10->12, 11->13), leaving a gap (rows 10-11 blank) so
# the two entries can be split out by project part.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move row 11 (Postmortem) down to row 13 first, then row 10 (Testing)
# down to row 12, so we never overwrite a source row before reading it.
$ws.Range("B11").Copy($ws.Range("B13"))
$ws.Range("C11").Copy($ws.Range("C13"))
$ws.Range("D11").Copy($ws.Range("D13"))
$ws.Range("F11").Copy($ws.Range("F13"))
$ws.Range("H11").Copy($ws.Range("H13"))

$ws.Range("B10").Copy($ws.Range("B12"))
$ws.Range("C10").Copy($ws.Range("C12"))
$ws.Range("D10").Copy($ws.Range("D12"))
$ws.Range("F10").Copy($ws.Range("F12"))
$ws.Range("H10").Copy($ws.Range("H12"))

# Clear out the old rows 10-11 entirely, leaving a blank gap.
$ws.Range("A10:H11").Clear()

# Match the new selection left behind in the saved file.
$ws.Range("C18").Select()
